$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 6 new rows right after row 6 (before the existing "K_FB_rauchkind_15" row)
$ws.Rows("7:12").Insert()

# The inserted rows inherit formatting from the row above; strip that so the
# new cells end up with the default (unstyled) formatting, matching the target.
$ws.Range("A7:D12").ClearFormats()

# New variable definitions (name / label / valueType) for the 6 inserted rows.
$names  = @("aktivleiisommn_15", "aktivleiwintn_15", "aktivmitsommn_15", "aktivmitwintn_15", "aktivschsommn_15", "aktivschwintn_15")
$labels = @(
    "Number of hours in a normal week (7 days) with light physical activity of the child (without sweating, normal breathing, e.g. walking) in summer [hours]",
    "Number of hours in a normal week (7 days) with light physical activity of the child (without sweating, normal breathing, e.g. walking) in winter [hours]",
    "Number of hours in a normal week (7 days) with moderate physical activity (a little sweating, slightly increased breathing e.g. cycling, swimming, skating) in summer [hours]",
    "Number of hours in a normal week (7 days) with moderate physical activity (some sweating, slightly increased breathing e.g. cycling, swimming, skating) in winter [hours]",
    "Number of hours in a normal week (7 days) with heavy physical activity (a lot of sweating, rapid breathing, e.g. ball games, training) in summer [hours]",
    "Number of hours in a normal week (7 days) with heavy physical activity (a lot of sweating, rapid breathing, e.g. ball games, training) in winter [hours]"
)

for ($i = 0; $i -lt 6; $i++) {
    $row = 7 + $i
    $ws.Range("B$row").Value2 = $names[$i]
    $ws.Range("C$row").Value2 = $labels[$i]
    $ws.Range("D$row").Value2 = "decimal"
}

# The "index" column (A) is a plain sequential count (row number - 1). Rewrite
# it for every data row from the first inserted row through the last row so
# the numbering stays consistent after the insert.
$lastRow = $ws.UsedRange.Rows.Count
$n = $lastRow - 7 + 1
$idx = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $idx[$i,0] = 6 + $i
}
$ws.Range("A7:A$lastRow").Value2 = $idx

# Mirror the saved selection state recorded in the workbook.
$ws.Range("A2:A$lastRow").Select()

$wb.Save()
